# Update recomputed NATMI TPM-derived metrics for the Tnfsf13 -> Tnfrsf1a LR-pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.6848073333333332
$ws.Range("H2").Value = 2.054422
$ws.Range("I2").Value = 0.2268310526442471
$ws.Range("J2").Value = 0.2268310526442472
$ws.Range("M2").Value = 14.65767833333333
$ws.Range("N2").Value = 43.973035
$ws.Range("O2").Value = 0.2345581433878666
$ws.Range("P2").Value = 0.2345581433878665
$ws.Range("Q2").Value = 10.03768561230778
$ws.Range("R2").Value = 90.33917051076999
$ws.Range("S2").Value = 0.05320507057095002
$ws.Range("T2").Value = 0.05320507057095003

# Row 3
$ws.Range("G3").Value = 0.6848073333333332
$ws.Range("H3").Value = 2.054422
$ws.Range("I3").Value = 0.2268310526442471
$ws.Range("J3").Value = 0.2268310526442472
$ws.Range("M3").Value = 31.695371
$ws.Range("N3").Value = 95.086113
$ws.Range("O3").Value = 0.5072022462686253
$ws.Range("P3").Value = 0.5072022462686253
$ws.Range("Q3").Value = 21.70522249352066
$ws.Range("R3").Value = 195.347002441686
$ws.Range("S3").Value = 0.1150492194246389
$ws.Range("T3").Value = 0.115049219424639

# Row 4
$ws.Range("G4").Value = 0.6848073333333332
$ws.Range("H4").Value = 2.054422
$ws.Range("I4").Value = 0.2268310526442471
$ws.Range("J4").Value = 0.2268310526442472
$ws.Range("M4").Value = 16.13754733333333
$ws.Range("N4").Value = 48.41264200000001
$ws.Range("O4").Value = 0.2582396103435082
$ws.Range("P4").Value = 0.2582396103435082
$ws.Range("Q4").Value = 11.05111075588044
$ws.Range("R4").Value = 99.459996802924
$ws.Range("S4").Value = 0.05857676264865817
$ws.Range("T4").Value = 0.05857676264865819

# Row 5
$ws.Range("I5").Value = 0.1086184939966157
$ws.Range("J5").Value = 0.1086184939966157
$ws.Range("M5").Value = 14.65767833333333
$ws.Range("N5").Value = 43.973035
$ws.Range("O5").Value = 0.2345581433878666
$ws.Range("P5").Value = 0.2345581433878665
$ws.Range("Q5").Value = 4.806565422637778
$ws.Range("R5").Value = 43.25908880374
$ws.Range("S5").Value = 0.0254773522894323
$ws.Range("T5").Value = 0.0254773522894323

# Row 6
$ws.Range("I6").Value = 0.1086184939966157
$ws.Range("J6").Value = 0.1086184939966157
$ws.Range("M6").Value = 31.695371
$ws.Range("N6").Value = 95.086113
$ws.Range("O6").Value = 0.5072022462686253
$ws.Range("P6").Value = 0.5072022462686253
$ws.Range("R6").Value = 93.542294869332
$ws.Range("S6").Value = 0.05509154414139866
$ws.Range("T6").Value = 0.05509154414139866

# Row 7
$ws.Range("I7").Value = 0.1086184939966157
$ws.Range("J7").Value = 0.1086184939966157
$ws.Range("M7").Value = 16.13754733333333
$ws.Range("N7").Value = 48.41264200000001
$ws.Range("O7").Value = 0.2582396103435082
$ws.Range("P7").Value = 0.2582396103435082
$ws.Range("Q7").Value = 5.291846038276445
$ws.Range("R7").Value = 47.626614344488
$ws.Range("S7").Value = 0.02804959756578472
$ws.Range("T7").Value = 0.02804959756578472

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.006290666666667
$ws.Range("H8").Value = 6.018872
$ws.Range("I8").Value = 0.6645504533591371
$ws.Range("J8").Value = 0.6645504533591372
$ws.Range("M8").Value = 14.65767833333333
$ws.Range("N8").Value = 43.973035
$ws.Range("O8").Value = 0.2345581433878666
$ws.Range("P8").Value = 0.2345581433878665
$ws.Range("Q8").Value = 29.40756323516889
$ws.Range("R8").Value = 264.66806911652
$ws.Range("S8").Value = 0.1558757205274842
$ws.Range("T8").Value = 0.1558757205274842

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.006290666666667
$ws.Range("H9").Value = 6.018872
$ws.Range("I9").Value = 0.6645504533591371
$ws.Range("J9").Value = 0.6645504533591372
$ws.Range("M9").Value = 31.695371
$ws.Range("N9").Value = 95.086113
$ws.Range("O9").Value = 0.5072022462686253
$ws.Range("P9").Value = 0.5072022462686253
$ws.Range("Q9").Value = 63.59012701383733
$ws.Range("R9").Value = 572.311143124536
$ws.Range("S9").Value = 0.3370614827025877
$ws.Range("T9").Value = 0.3370614827025877

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.006290666666667
$ws.Range("H10").Value = 6.018872
$ws.Range("I10").Value = 0.6645504533591371
$ws.Range("J10").Value = 0.6645504533591372
$ws.Range("M10").Value = 16.13754733333333
$ws.Range("N10").Value = 48.41264200000001
$ws.Range("O10").Value = 0.2582396103435082
$ws.Range("P10").Value = 0.2582396103435082
$ws.Range("Q10").Value = 32.37661059775822
$ws.Range("R10").Value = 291.389495379824
$ws.Range("S10").Value = 0.1716132501290653
$ws.Range("T10").Value = 0.1716132501290653
